$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the code text for A3 (merge manali's code: "print 'Hey';" -> "print 'Hey Selenium';")
$ws.Range("A3").Value = "print 'Hey Selenium';"

# Update the selected cell to A3 to match the saved selection state
$ws.Range("A3").Select()
